$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.518.30"
$ws.Range("D3").Value = "2.098.63"
$ws.Range("E3").Value = "  +9.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'252.60"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'0.655"
$ws.Range("E6").Value = "  -6.68%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'47.36"
$ws.Range("E8").Value = "  +3.55%  "
$ws.Range("D9").Value = "'60.35"
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("D12").Value = "'0.1000"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'14.56"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "2.402.59"
$ws.Range("E14").Value = "  +9.60%  "
$ws.Range("D15").Value = "'0.829"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "2.090.93"
$ws.Range("E16").Value = "  +9.23%  "
$ws.Range("D17").Value = "'5.08"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "36.476.87"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'72.75"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").Value = "0.0₃0829"
$ws.Range("E20").Value = "  -3.71%  "
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "'239.87"
$ws.Range("E22").Value = "  -4.07%  "
$ws.Range("D23").Value = "'5.18"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  -5.32%  "
$ws.Range("D26").Value = "'170.39"
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("D27").Value = "'21.40"
$ws.Range("E27").Value = "  +14.50%  "
$ws.Range("D28").Value = "'9.12"
$ws.Range("E28").Value = "  +4.06%  "
$ws.Range("E29").Value = "  -9.53%  "
$ws.Range("D30").Value = "'28.88"
$ws.Range("E30").Value = "  +62.57%  "
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("D33").Value = "'0.0616"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'0.993"
$ws.Range("E34").Value = "  +13.51%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'2.43"
$ws.Range("E35").Value = "  +21.29%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.0911"
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").Value = "'4.10"
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("E40").Value = "  -11.40%  "
$ws.Range("E41").Value = "  +6.30%  "
$ws.Range("D42").Value = "'0.0223"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "'97.51"
$ws.Range("E43").Value = "  -7.27%  "
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").Value = "'15.87"
$ws.Range("E45").Value = "  -8.84%  "
$ws.Range("D46").Value = "1.326.90"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "'0.0843"
$ws.Range("E47").Value = "  +3.62%  "
$ws.Range("D48").Value = "'7.05"
$ws.Range("E48").Value = "  +9.29%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.84"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.284.91"
$ws.Range("E50").Value = "  +9.43%  "
$ws.Range("D51").Value = "'2.24"
$ws.Range("E51").Value = "  -5.81%  "
